$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Centrophorus granulosus / Gulper shark
$ws.Range("C2").Value = "Not applicable"
$ws.Range("D2").Value = "Not applicable"
$ws.Range("E2").Value = "Not applicable"
$ws.Range("F2").Value = "Unknown"
$ws.Range("G2").Value = "Unknown"

# Row 3 - Centrophorus squamosus / Leafscale gulper shark
$ws.Range("C3").Value = "Unknown"
$ws.Range("D3").Value = "Not applicable"

# Row 4 - Centroscymnus coelolepis / Portuguese dogfish
$ws.Range("D4").Value = "Not applicable"

# Row 6 - Dipturus batis / Blue skate
$ws.Range("E6").Value = "Unknown"
$ws.Range("F6").Value = "Unknown"
$ws.Range("G6").Value = "Unknown"

# Row 7 - Lamna nasus / Porbeagle
$ws.Range("C7").Value = "Unknown"
$ws.Range("D7").Value = "Unknown"
$ws.Range("E7").Value = "Unknown"
$ws.Range("F7").Value = "Unknown"
$ws.Range("G7").Value = "Unknown"

# Row 8 - Raja clavata / Thornback ray
$ws.Range("C8").Value = "Unknown"
$ws.Range("G8").Value = "Stable"

# Row 9 - Raja montagui / Spotted ray
$ws.Range("C9").Value = "Not applicable"
$ws.Range("E9").Value = "Unknown"
$ws.Range("F9").Value = "Unknown"
$ws.Range("G9").Value = "Not applicable"

# Row 10 - Rostroraja alba / White skate
$ws.Range("C10").Value = "Not applicable"
$ws.Range("G10").Value = "Not applicable"

# Row 11 - Squatina squatina / Angelshark
$ws.Range("C11").Value = "Not applicable"
$ws.Range("E11").Value = "Poor"
$ws.Range("G11").Value = "Not applicable"
